$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet lists reactions grouped by ligand, each group being a bold
# header row (ligand name) followed by its reaction rows. We are adding
# two new groups - "VEGF-B" and "PlGF" - right before the existing
# "PDGF-AA" group (i.e. after row 9 / before row 10), each group being a
# header row + 2 reaction rows, for a total of 6 new rows.

# 1) Insert 6 blank rows at row 10, shifting everything at/after row 10
#    down by 6 rows (PDGF-AA group, formerly at rows 10-12, moves to
#    rows 16-18, etc).
$ws.Range("A10:C15").EntireRow.Insert()

# 2) Populate formatting for the 6 new (currently blank) rows by copying
#    the style of the PDGF-AA group (now shifted to rows 16-18), which
#    has the exact same shape as our new groups: 1 header row + 2
#    reaction rows.
$ws.Range("A16:C18").Copy($ws.Range("A10"))
$ws.Range("A16:C18").Copy($ws.Range("A13"))

# 3) Fill in the text for the new "VEGF-B" group (rows 10-12).
$ws.Range("A10").Value = "VEGF-B"

$ws.Range("A11").Value = "VB + R1 <-> VB:R1"
$ws.Range("B11").Value = "konVBR1"
$ws.Range("C11").Value = "koffVBR1"

$ws.Range("A12").Value = "VB + N1 <-> VB:N1"
$ws.Range("B12").Value = "konVBN1"
$ws.Range("C12").Value = "koffVBN1"

# 4) Fill in the text for the new "PlGF" group (rows 13-15).
$ws.Range("A13").Value = "PlGF"

$ws.Range("A14").Value = "Pl + R1 <-> Pl:R1"
$ws.Range("B14").Value = "konPlR1"
$ws.Range("C14").Value = "koffPlR1"

$ws.Range("A15").Value = "Pl + N1 <-> Pl:N1"
$ws.Range("B15").Value = "konPlN1"
$ws.Range("C15").Value = "koffPlN1"

# 5) Match the saved selection state.
[void]$ws.Range("B13").Select()
